# Regenerate s_vals data (filter save games) for merryweather_julian.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TB(B), d2S(C), K(D), IP(E) and recomputed sum(G) per row (rows 2-21)
$rows = @(
  @{ Row=2; B=3.182878228561681; C=1.65323645889881; D=0.7127328510149897; E=0.4998867070740569; G=6.048734245549538 },
  @{ Row=3; B=0.02258322285507441; C=1.65323645889881; D=0.1529057820181812; E=0.4998867070740569; G=2.328612170846122 },
  @{ Row=4; B=3.182878228561681; C=1.65323645889881; D=3.082599426703578; E=6.48142807727062; G=14.40014219143469 },
  @{ Row=5; B=3.182878228561681; C=1.65323645889881; D=3.082599426703578; E=0.4998867070740569; G=8.418600821238126 },
  @{ Row=6; B=0.1554434735375247; C=0.3375848360084654; D=0.7127328510149897; E=0.4998867070740569; G=1.705647867635037 },
  @{ Row=7; B=3.182878228561681; C=1.65323645889881; D=3.082599426703578; E=0.4998867070740569; G=8.418600821238126 },
  @{ Row=8; B=3.182878228561681; C=1.65323645889881; D=3.082599426703578; E=0.4998867070740569; G=8.418600821238126 },
  @{ Row=9; B=3.182878228561681; C=1.65323645889881; D=0.7127328510149897; E=0.4998867070740569; G=6.048734245549538 },
  @{ Row=10; B=3.182878228561681; C=1.65323645889881; D=3.082599426703578; E=0.4998867070740569; G=8.418600821238126 },
  @{ Row=11; B=3.182878228561681; C=1.65323645889881; D=0.1529057820181812; E=0.4998867070740569; G=5.488907176552729 },
  @{ Row=12; B=0.7287194209349384; C=0.3375848360084654; D=3.082599426703578; E=0.4998867070740569; G=4.64879039072104 },
  @{ Row=13; B=3.182878228561681; C=1.65323645889881; D=0.7127328510149897; E=0.4998867070740569; G=6.048734245549538 },
  @{ Row=14; B=3.182878228561681; C=1.65323645889881; D=0.1529057820181812; E=0.4998867070740569; G=5.488907176552729 },
  @{ Row=15; B=1.505614041169197; C=1.65323645889881; D=3.082599426703578; E=0.4998867070740569; G=6.741336633845642 },
  @{ Row=16; B=1.505614041169197; C=1.65323645889881; D=0.7127328510149897; E=0.4998867070740569; G=4.371470058157054 },
  @{ Row=17; B=0.3464964993005633; C=0.3375848360084654; D=3.082599426703578; E=6.48142807727062; G=10.24810883928323 },
  @{ Row=18; B=3.182878228561681; C=1.65323645889881; D=0.1529057820181812; E=0.4998867070740569; G=5.488907176552729 },
  @{ Row=19; B=1.505614041169197; C=1.65323645889881; D=0.7127328510149897; E=0.4998867070740569; G=4.371470058157054 },
  @{ Row=20; B=3.182878228561681; C=1.65323645889881; D=3.082599426703578; E=0.4998867070740569; G=8.418600821238126 },
  @{ Row=21; B=0.7287194209349384; C=0.3375848360084654; D=3.082599426703578; E=0.4998867070740569; G=4.64879039072104 }
)

foreach ($r in $rows) {
  $ws.Range("B$($r.Row)").Value = $r.B
  $ws.Range("C$($r.Row)").Value = $r.C
  $ws.Range("D$($r.Row)").Value = $r.D
  $ws.Range("E$($r.Row)").Value = $r.E
  $ws.Range("G$($r.Row)").Value = $r.G
}
